# R_act_qualifier_cockpit.xlsx - "actual qualifier ranks | per ctr cockpit"
#
# The header row (row 1) holds one label per data column (V1..V160, one per
# qualifying-position column). This edit rotates the header labels one slot
# to the right: every label moves from column c to column c+1, and the label
# that falls off the end (the old last column, "V160") wraps around to
# become the new first column - renamed from "V160" to "V0" (the "actual"/
# current rank column called out in the commit message).
#
# Net effect matches the diff exactly:
#   A1 <- old last column's label, retitled "V0"
#   B1 <- old A1's label
#   C1 <- old B1's label
#   ...
#   (last col) <- old second-to-last column's label

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$colCount = $usedRange.Columns.Count

# Snapshot the current header row values before overwriting anything
# (.Value2 reads the real scalar; .Value on this host returns a property
# descriptor string instead of the cell's contents).
$headers = @()
for ($c = 1; $c -le $colCount; $c++) {
    $headers += $ws.Cells.Item(1, $c).Value2
}

# Shift every label one column to the right (B1=old A1, C1=old B1, ...).
for ($c = $colCount; $c -ge 2; $c--) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 2]
}

# The old last column's label wraps to the front, renamed V160 -> V0.
$ws.Cells.Item(1, 1).Value = "V0"

# Match the author's final selection.
$null = $ws.Range("FD1").Select()
